$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.293.55'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '3.670.17'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '645.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.497'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.06'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.443'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000230'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").Value = '4.289.11'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").Value = '3.652.45'
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '69.306.71'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.116'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.643'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = '3.817.07'
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.32%  '
$ws.Range("E29").Value = '  -3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.66%  '
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.662.63'
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.162'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '178.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.90%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.74%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0888'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.923'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.71%  '
$ws.Range("E48").Value = '  -4.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000264'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.09%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.10%  '
